$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Corrections to existing entries (end-time edits on rows 4, 5, 11) ---
$ws.Range("E4").Value = 0.46527777777777773
$ws.Range("E5").Value = 0.60069444444444442
$ws.Range("E11").Value = 0.41666666666666669

# --- New time-registration entries in rows 15-19 ---
# Row 15: OC0802 rettelser / System Analyst
$ws.Range("A15").Value = "OC0802 rettelser"
$ws.Range("B15").Value = "System Analyst "
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 43889
$ws.Range("D15").Value = 0.36458333333333331
$ws.Range("E15").Value = 0.38541666666666669

# Row 17: ATD07 (shared string inserted before "ATD06" to match authoring order)
$ws.Range("A17").Value = "ATD07"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = 43889
$ws.Range("D17").Value = 0.42708333333333331
$ws.Range("E17").Value = 0.45833333333333331

# Row 16: ATD06
$ws.Range("A16").Value = "ATD06"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = 43889
$ws.Range("D16").Value = 0.3888888888888889
$ws.Range("E16").Value = 0.42708333333333331

# Row 18: ATD07 rettelser
$ws.Range("A18").Value = "ATD07 rettelser"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = 43889
$ws.Range("D18").Value = 0.49305555555555558
$ws.Range("E18").Value = 0.52430555555555558

# Row 19: ATD06 rettelser
$ws.Range("A19").Value = "ATD06 rettelser"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value = 43889
$ws.Range("D19").Value = 0.52638888888888891
$ws.Range("E19").Value = 0.56944444444444442

# --- Update the active selection to match the author's last position ---
$ws.Range("B20").Select() | Out-Null
